# Generate Report for Handoff
# The 69f9321c-... item has moved from "In Translation" to "Ready for
# handoff" and is now listed after b93bf6ae-... in each report sheet
# (Overview, zh-cn, de-de). Swap the two data rows on every sheet and
# update the 69f9321c row's status / handoff timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "b93bf6ae-3694-4483-95c5-446da0edd3f1.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "2016-03-23 02:16:51"

$ws.Range("A3").Value = "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-23 02:17:47"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2096270b40e1f0937a7ccda4783bb908e23ce736/e2e/69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md", "", "", "b93bf6ae-3694-4483-95c5-446da0edd3f1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2096270b40e1f0937a7ccda4783bb908e23ce736/e2e/b93bf6ae-3694-4483-95c5-446da0edd3f1.md", "", "", "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "b93bf6ae-3694-4483-95c5-446da0edd3f1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "b93bf6ae-3694-4483-95c5-446da0edd3f1.5f16319619eea5a3a5ce01852e4adf0e2af98a73.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-23 02:16:47"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-23 02:17:42"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2096270b40e1f0937a7ccda4783bb908e23ce736/e2e/69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md", "", "", "b93bf6ae-3694-4483-95c5-446da0edd3f1.md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c178f567366241603952672e1cc2070a5ac1015/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.zh-cn.xlf", "", "", "b93bf6ae-3694-4483-95c5-446da0edd3f1.5f16319619eea5a3a5ce01852e4adf0e2af98a73.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2096270b40e1f0937a7ccda4783bb908e23ce736/e2e/b93bf6ae-3694-4483-95c5-446da0edd3f1.md", "", "", "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c178f567366241603952672e1cc2070a5ac1015/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b93bf6ae-3694-4483-95c5-446da0edd3f1.5f16319619eea5a3a5ce01852e4adf0e2af98a73.zh-cn.xlf", "", "", "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "b93bf6ae-3694-4483-95c5-446da0edd3f1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "b93bf6ae-3694-4483-95c5-446da0edd3f1.5f16319619eea5a3a5ce01852e4adf0e2af98a73.de-de.xlf"
$ws.Range("E2").Value = "2016-03-23 02:16:51"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.de-de.xlf"
$ws.Range("E3").Value = "2016-03-23 02:17:47"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2096270b40e1f0937a7ccda4783bb908e23ce736/e2e/69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md", "", "", "b93bf6ae-3694-4483-95c5-446da0edd3f1.md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/812a6328a1d6faf719c6b2e79b6cda4b8b9276ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.de-de.xlf", "", "", "b93bf6ae-3694-4483-95c5-446da0edd3f1.5f16319619eea5a3a5ce01852e4adf0e2af98a73.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2096270b40e1f0937a7ccda4783bb908e23ce736/e2e/b93bf6ae-3694-4483-95c5-446da0edd3f1.md", "", "", "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/812a6328a1d6faf719c6b2e79b6cda4b8b9276ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b93bf6ae-3694-4483-95c5-446da0edd3f1.5f16319619eea5a3a5ce01852e4adf0e2af98a73.de-de.xlf", "", "", "69f9321c-0d3f-49c0-9e04-6b0ac2a1a391.229886f1f51574ec483b3bf868b13ff072de2820.de-de.xlf")

Write-Output "Report regenerated for handoff"
